# Il12a-Il12rb1.xlsx — refresh with new TPM-derived values.
#
# The "ECs" sending-cluster block (old rows 2-4) is dropped entirely, and the
# remaining Sending-cluster=FAPs / Sending-cluster=MuSCs blocks (old rows
# 5-10) move up to become rows 2-7, with every numeric column recomputed
# against the updated TPM data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the old "ECs" sending-cluster rows; rows below shift up automatically.
$ws.Rows("2:4").Delete()

# 2) Write the refreshed rows (A:D = labels, E:T = new TPM-derived numbers).
$data = @(
    @(2, "FAPs",  "Il12a", "Il12rb1", "ECs",   3, 1,                  1.574655333333333, 4.723966, 0.805795973174511,  0.8057959731745109, 2, 0.6666666666666666, 0.07429999999999999, 0.2229,   0.08480128194885443, 0.08480128194885443, 0.1169968912666667,  1.0529720214,       0.06833253151442326, 0.06833253151442324),
    @(3, "FAPs",  "Il12a", "Il12rb1", "FAPs",  3, 1,                  1.574655333333333, 4.723966, 0.805795973174511,  0.8057959731745109, 3, 1,                  0.6279076666666666,  1.883723, 0.7166537695672586,  0.7166537695672586,  0.9887381561575554,  8.898643405417999,  0.5774767216776309, 0.5774767216776308),
    @(4, "FAPs",  "Il12a", "Il12rb1", "MuSCs", 3, 1,                  1.574655333333333, 4.723966, 0.805795973174511,  0.8057959731745109, 2, 0.6666666666666666, 0.1739583333333333,  0.521875, 0.198544948483887,   0.198544948483887,   0.2739244173611111,  2.46531975625,      0.1599867199824569, 0.1599867199824568),
    @(5, "MuSCs", "Il12a", "Il12rb1", "ECs",   2, 0.6666666666666666, 0.379506,           1.138518, 0.194204026825489,  0.194204026825489,  2, 0.6666666666666666, 0.07429999999999999, 0.2229,   0.08480128194885443, 0.08480128194885443, 0.02819729579999999, 0.2537756622,       0.01646875043443118, 0.01646875043443118),
    @(6, "MuSCs", "Il12a", "Il12rb1", "FAPs",  2, 0.6666666666666666, 0.379506,           1.138518, 0.194204026825489,  0.194204026825489,  3, 1,                  0.6279076666666666,  1.883723, 0.7166537695672586,  0.7166537695672586,  0.238294726946,      2.144652542514,     0.1391770478896277, 0.1391770478896277),
    @(7, "MuSCs", "Il12a", "Il12rb1", "MuSCs", 2, 0.6666666666666666, 0.379506,           1.138518, 0.194204026825489,  0.194204026825489,  2, 0.6666666666666666, 0.1739583333333333,  0.521875, 0.198544948483887,   0.198544948483887,   0.06601823124999999, 0.59416408125,      0.03855822850143012, 0.03855822850143011)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 1; $i -lt $row.Length; $i++) {
        $col = $cols[$i - 1]
        $ws.Range("$col$r").Value = $row[$i]
    }
}
